# Add the new "Induction" diagram slide (slide 9 / sldId 306) at the end
# of the deck, reproducing the pyramid-of-connectors + circle + dashed
# leader-line figure described by the target diff.

$p = $ppt.ActivePresentation

# --- 1. Insert a new blank slide at the end -----------------------------
$newIndex = $p.Slides.Count + 1
$slide = $p.Slides.Add($newIndex, 12)   # ppLayoutBlank

# EMU -> point conversion (PowerPoint's object model works in points)
$emu = 12700.0

function ToPt([double]$v) { return $v / $emu }

# --- 2. The twelve "Straight Arrow Connector" diagonal lines ------------
# Each is a straight connector with a triangle arrowhead, 3pt weight,
# colour = Background 1, Darker 15% (bg1 / lumMod 85%).
$arrowData = @(
    @{ Name = "Straight Arrow Connector 4";  X = 1983608; Y = 606390; CX = 962526; CY = 4263991 },
    @{ Name = "Straight Arrow Connector 5";  X = 2239881; Y = 606390; CX = 962526; CY = 4263991 },
    @{ Name = "Straight Arrow Connector 6";  X = 2493344; Y = 606390; CX = 962526; CY = 4263991 },
    @{ Name = "Straight Arrow Connector 7";  X = 2760044; Y = 606391; CX = 962526; CY = 4263991 },
    @{ Name = "Straight Arrow Connector 8";  X = 3039176; Y = 606391; CX = 962526; CY = 4263991 },
    @{ Name = "Straight Arrow Connector 9";  X = 3317507; Y = 606392; CX = 962526; CY = 4263991 },
    @{ Name = "Straight Arrow Connector 10"; X = 3599048; Y = 606388; CX = 962526; CY = 4263991 },
    @{ Name = "Straight Arrow Connector 11"; X = 3855321; Y = 606388; CX = 962526; CY = 4263991 },
    @{ Name = "Straight Arrow Connector 12"; X = 4108784; Y = 606388; CX = 962526; CY = 4263991 },
    @{ Name = "Straight Arrow Connector 13"; X = 4375484; Y = 606389; CX = 962526; CY = 4263991 },
    @{ Name = "Straight Arrow Connector 14"; X = 4654616; Y = 606389; CX = 962526; CY = 4263991 },
    @{ Name = "Straight Arrow Connector 15"; X = 4932947; Y = 606390; CX = 962526; CY = 4263991 }
)

foreach ($d in $arrowData) {
    $x1 = ToPt($d.X)
    $y1 = ToPt($d.Y)
    $x2 = ToPt($d.X + $d.CX)
    $y2 = ToPt($d.Y + $d.CY)

    $cxn = $slide.Shapes.AddConnector(1, $x1, $y1, $x2, $y2)   # msoConnectorStraight
    $cxn.Name = $d.Name
    $cxn.Line.Weight = 3
    $cxn.Line.ForeColor.RGB = 0xD9D9D9
    $cxn.Line.EndArrowheadStyle = 2   # msoArrowheadTriangle
}

# --- 3. The horizontal "Straight Arrow Connector 17" (id 18) ------------
$hx1 = ToPt(3797165)
$hy1 = ToPt(2666198)
$hx2 = ToPt(3797165 + 2203385)
$hy2 = ToPt(2666198 + 0)
$hconn = $slide.Shapes.AddConnector(1, $hx1, $hy1, $hx2, $hy2)
$hconn.Name = "Straight Arrow Connector 17"
$hconn.Line.Weight = 3
$hconn.Line.ForeColor.RGB = 0x000000
$hconn.Line.EndArrowheadStyle = 2

# --- 4. The filled black circle ("Oval 18") ------------------------------
$ox = ToPt(3608269)
$oy = ToPt(2503771)
$ow = ToPt(324853)
$oh = ToPt(324853)
$oval = $slide.Shapes.AddShape(9, $ox, $oy, $ow, $oh)   # msoShapeOval
$oval.Name = "Oval 18"
$oval.Fill.ForeColor.RGB = 0x000000
$oval.Line.ForeColor.RGB = 0x000000
$oval.Line.Weight = 3
$oval.TextFrame.TextRange.ParagraphFormat.Alignment = 2   # ppAlignCenter

# --- 5. The dashed leader line ("Straight Connector 20") ----------------
$lx1 = ToPt(3722570)
$ly1 = ToPt(433137)
$lx2 = ToPt(3722570 + 48125)
$ly2 = ToPt(433137 + 2233060)
$dash = $slide.Shapes.AddLine($lx1, $ly1, $lx2, $ly2)
$dash.Name = "Straight Connector 20"
$dash.HorizontalFlip = -1
$dash.VerticalFlip = -1
$dash.Line.Weight = 3
$dash.Line.ForeColor.RGB = 0x000000
$dash.Line.DashStyle = 4   # msoLineDash
$dash.Line.EndArrowheadStyle = 2
